$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.433.61'
$ws.Range('E2').Value = '  -1.06%  '
$ws.Range('D3').Value = '2.406.73'
$ws.Range('E3').Value = '  -1.90%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '507.33'
$ws.Range('E5').Value = '  -3.37%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '132.95'
$ws.Range('E6').Value = '  +0.97%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.995'
$ws.Range('E7').Value = '  -0.51%  '
$ws.Range('E8').Value = '  -1.30%  '
$ws.Range('D9').Value = '2.444.12'
$ws.Range('E9').Value = '  -0.55%  '
$ws.Range('E10').Value = '  -0.28%  '
$ws.Range('E11').Value = '  -1.32%  '
$ws.Range('E12').Value = '  -0.83%  '
$ws.Range('E13').Value = '  -7.71%  '
$ws.Range('D14').Value = '2.842.65'
$ws.Range('E14').Value = '  -1.54%  '
$ws.Range('D15').Value = '57.282.37'
$ws.Range('E15').Value = '  -1.12%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '21.85'
$ws.Range('E16').Value = '  +0.33%  '
$ws.Range('E17').Value = '  +0.12%  '
$ws.Range('D18').Value = '2.463.65'
$ws.Range('E18').Value = '  +0.29%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '10.28'
$ws.Range('E19').Value = '  -0.67%  '
$ws.Range('E20').Value = '  -0.33%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '313.78'
$ws.Range('E21').Value = '  +0.41%  '
$ws.Range('E22').Value = '  +5.18%  '
$ws.Range('E23').Value = '  -0.30%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '5.70'
$ws.Range('E24').Value = '  -2.55%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '65.11'
$ws.Range('E25').Value = '  +0.18%  '
$ws.Range('E26').Value = '  -0.69%  '
$ws.Range('D27').Value = '2.524.59'
$ws.Range('E27').Value = '  -2.77%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.381'
$ws.Range('E28').Value = '  -5.83%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.154'
$ws.Range('E29').Value = '  -2.32%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '7.56'
$ws.Range('E30').Value = '  +4.39%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '173.74'
$ws.Range('E31').Value = '  +0.20%  '
$ws.Range('E32').Value = '  -0.88%  '
$ws.Range('E33').Value = '  -1.00%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '6.16'
$ws.Range('E34').Value = '  -1.38%  '
$ws.Range('E35').Value = '  -1.30%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.992'
$ws.Range('E37').Value = '  -0.44%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '17.96'
$ws.Range('E38').Value = '  +0.85%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.23'
$ws.Range('E39').Value = '  +3.45%  '
$ws.Range('E40').Value = '  +0.52%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '36.64'
$ws.Range('E41').Value = '  +1.04%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.813'
$ws.Range('E42').Value = '  -0.24%  '
$ws.Range('E43').Value = '  +0.36%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '135.12'
$ws.Range('E44').Value = '  +10.40%  '
$ws.Range('B45').Value = 'Filecoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.39'
$ws.Range('E45').Value = '  -0.69%  '
$ws.Range('B46').Value = 'RenderToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '5.00'
$ws.Range('E46').Value = '  +3.75%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '255.13'
$ws.Range('E47').Value = '  -2.91%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.571'
$ws.Range('E48').Value = '  -2.82%  '
$ws.Range('E49').Value = '  -0.60%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0492'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0213'
$ws.Range('E51').Value = '  +0.47%  '
